$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9631166458129883
$ws.Range("B1").Value = 2.226003646850586
$ws.Range("C1").Value = 8.196019172668457
$ws.Range("D1").Value = 1.803224921226501
$ws.Range("E1").Value = 1.289706468582153
